# Update the "dSF" (column F) values for the rows whose data was re-pulled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -3
    8  = -2
    9  = 1
    10 = 3
    11 = 3
    25 = -1
    26 = 1
    38 = 2
    39 = 2
    48 = -1
    52 = 0
    53 = 1
    55 = -1
    56 = -1
    59 = -1
    66 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
